$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 399
$ws.Range("I6").Value = 399
$ws.Range("K6").Value = 1197
$ws.Range("M6").Value = -1085

$ws.Range("H11").Value = 153.44444
$ws.Range("I11").Value = 153.44444
$ws.Range("K11").Value = 153.44444
$ws.Range("M11").Value = -13.44443999999999

$ws.Range("H137").Value = 2871.5862
$ws.Range("J137").Value = 4605.8667
$ws.Range("L137").Value = 13817.6001
$ws.Range("N137").Value = -18917.6001

$ws.Range("H138").Value = 5737.4116
$ws.Range("I138").Value = 1378.7142
$ws.Range("J138").Value = 8788.5
$ws.Range("K138").Value = 4136.142599999999
$ws.Range("L138").Value = 26365.5
$ws.Range("M138").Value = 1003.857400000001
$ws.Range("N138").Value = -36645.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("K4").Value = 100
$ws.Range("M4").Value = 16

$ws.Range("H61").Value = 3784.5715
$ws.Range("I61").Value = 3845.1538
$ws.Range("K61").Value = 3845.1538
$ws.Range("M61").Value = -3633.1538

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("N110").ClearContents()

$ws.Range("H122").Value = 907.8
$ws.Range("I122").Value = 876.6842
$ws.Range("K122").Value = 2630.0526
$ws.Range("M122").Value = -180.0526

$ws.Range("H132").Value = 1771.3235
$ws.Range("I132").Value = 1455.9
$ws.Range("K132").Value = 4367.700000000001
$ws.Range("M132").Value = -1837.700000000001

$ws.Range("H136").Value = 3784.5715
$ws.Range("I136").Value = 3845.1538
$ws.Range("K136").Value = 11535.4614
$ws.Range("M136").Value = -8985.4614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3186.1
$ws.Range("I20").Value = 3216.5715
$ws.Range("K20").Value = 3216.5715
$ws.Range("M20").Value = -2969.5715

$ws.Range("H64").Value = 7
$ws.Range("J64").Value = 7
$ws.Range("L64").Value = 7
$ws.Range("N64").Value = -457

$ws.Range("H67").Value = 7
$ws.Range("J67").Value = 7
$ws.Range("L67").Value = 7
$ws.Range("N67").Value = -1567

$ws.Range("H107").Value = 642.2
$ws.Range("I107").Value = 728
$ws.Range("K107").Value = 728
$ws.Range("M107").Value = 1192

$ws.Range("H134").Value = 3786.1875
$ws.Range("I134").Value = 3272.5386
$ws.Range("K134").Value = 9817.6158
$ws.Range("M134").Value = -7282.6158

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1838.3334
$ws.Range("I31").Value = 1857.6
$ws.Range("J31").Value = 1799.8
$ws.Range("K31").Value = 1857.6
$ws.Range("L31").Value = 1799.8
$ws.Range("M31").Value = -1562.6
$ws.Range("N31").Value = -2389.8

$ws.Range("H34").Value = 1838.3334
$ws.Range("I34").Value = 1857.6
$ws.Range("J34").Value = 1799.8
$ws.Range("K34").Value = 1857.6
$ws.Range("L34").Value = 1799.8
$ws.Range("M34").Value = -1655.6
$ws.Range("N34").Value = -2203.8

$ws.Range("H58").Value = 2442.5454
$ws.Range("I58").Value = 2298.5
$ws.Range("J58").Value = 2524.8572
$ws.Range("K58").Value = 2298.5
$ws.Range("L58").Value = 2524.8572
$ws.Range("M58").Value = -2095.5
$ws.Range("N58").Value = -2930.8572

$ws.Range("H62").Value = 4128.75
$ws.Range("I62").Value = 4087.6
$ws.Range("J62").Value = 4197.3335
$ws.Range("K62").Value = 4087.6
$ws.Range("L62").Value = 4197.3335
$ws.Range("M62").Value = -3463.6
$ws.Range("N62").Value = -5445.3335

$ws.Range("H65").Value = 4128.75
$ws.Range("I65").Value = 4087.6
$ws.Range("J65").Value = 4197.3335
$ws.Range("K65").Value = 20438
$ws.Range("L65").Value = 20986.6675
$ws.Range("M65").Value = -17318
$ws.Range("N65").Value = -27226.6675

$ws.Range("H132").Value = 2485.963
$ws.Range("I132").Value = 2192.8333
$ws.Range("K132").Value = 6578.499899999999
$ws.Range("M132").Value = -4048.499899999999

$ws.Range("H134").Value = 3526.2273
$ws.Range("I134").Value = 3509.5789
$ws.Range("K134").Value = 10528.7367
$ws.Range("M134").Value = -7993.736699999999

$ws.Range("H136").Value = 2442.5454
$ws.Range("I136").Value = 2298.5
$ws.Range("J136").Value = 2524.8572
$ws.Range("K136").Value = 6895.5
$ws.Range("L136").Value = 7574.571599999999
$ws.Range("M136").Value = -4345.5
$ws.Range("N136").Value = -12674.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 320.4
$ws.Range("I12").Value = 245
$ws.Range("J12").Value = 433.5
$ws.Range("K12").Value = 735
$ws.Range("L12").Value = 1300.5
$ws.Range("M12").Value = -562
$ws.Range("N12").Value = -1646.5

$ws.Range("H113").Value = 2795.9
$ws.Range("I113").Value = 998.5
$ws.Range("J113").Value = 3245.25
$ws.Range("K113").Value = 2995.5
$ws.Range("L113").Value = 9735.75
$ws.Range("M113").Value = -825.5
$ws.Range("N113").Value = -14075.75

$ws.Range("H122").Value = 1632.3334
$ws.Range("J122").Value = 2048.5
$ws.Range("L122").Value = 18436.5
$ws.Range("N122").Value = -23336.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4042.625
$ws.Range("I102").Value = 4042.625
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 4042.625
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -2420.625
$ws.Range("N102").ClearContents()

$ws.Range("H122").Value = 949
$ws.Range("I122").Value = 948.5
$ws.Range("K122").Value = 2845.5
$ws.Range("M122").Value = -395.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1100
$ws.Range("I22").Value = 1100
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1100
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -805
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 1100
$ws.Range("I27").Value = 1100
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1100
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -993
$ws.Range("N27").ClearContents()

$ws.Range("H46").Value = 4848.7
$ws.Range("J46").Value = 6414.5
$ws.Range("L46").Value = 6414.5
$ws.Range("N46").Value = -6790.5

$ws.Range("H55").Value = 213.18182
$ws.Range("I55").Value = 158.5
$ws.Range("J55").Value = 278.8
$ws.Range("K55").Value = 158.5
$ws.Range("L55").Value = 278.8
$ws.Range("M55").Value = 14.5
$ws.Range("N55").Value = -624.8

$ws.Range("H82").Value = 1408.6
$ws.Range("I82").Value = 1414.5
$ws.Range("J82").Value = 1399.75
$ws.Range("K82").Value = 1414.5
$ws.Range("L82").Value = 1399.75
$ws.Range("M82").Value = -1053.5
$ws.Range("N82").Value = -2121.75

$ws.Range("H85").Value = 1408.6
$ws.Range("I85").Value = 1414.5
$ws.Range("J85").Value = 1399.75
$ws.Range("K85").Value = 1414.5
$ws.Range("L85").Value = 1399.75
$ws.Range("M85").Value = -166.5
$ws.Range("N85").Value = -3895.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 20000
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws.Range("H62").Value = 13424.5
$ws.Range("I62").Value = 9899.666999999999
$ws.Range("K62").Value = 9899.666999999999
$ws.Range("M62").Value = -9275.666999999999

$ws.Range("H65").Value = 13424.5
$ws.Range("I65").Value = 9899.666999999999
$ws.Range("K65").Value = 49498.335
$ws.Range("M65").Value = -46378.335

$ws.Range("H132").Value = 2196.7812
$ws.Range("I132").Value = 1462.619
$ws.Range("K132").Value = 4387.857
$ws.Range("M132").Value = -1857.857
